$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Tyler Anderson: update Last.Updated date and Injury.Details text
$ws.Range("C2").Value = "July 21 2017"
$ws.Range("E2").Value = "Anderson has been designated for the 10-day disabled list as he will have arthroscopic left knee surgery and is expected to be sidelined until the beginning of August."

# Row 3 - Chad Bettis: update Last.Updated date and Injury.Details text
$ws.Range("C3").Value = "July 21 2017"
$ws.Range("E3").Value = "Bettis is on the 60-day disabled list while recovering from testicular cancer but is expected to make his season debut sometime during the early part of August."

# Row 4 - Tyler Chatwood: unchanged

# Row 5 - replace David Dahl entry with Gerardo Parra entry
$ws.Range("A5").Value = "Gerardo Parra"
$ws.Range("B5").Value = "parrage01"
$ws.Range("C5").Value = "July 23 2017"
$ws.Range("D5").Value = "Calf"
$ws.Range("E5").Value = "Parra sat out the previous game due to a right calf contusion and it is unclear as to if he will play against the Cardinals on Monday."

# Apply wrap text + explicit black font color to the Name column (A2:A5)
$ws.Range("A2:A5").WrapText = $true
$ws.Range("A2:A5").Font.Color = 0

# Update the saved selection to A23
$ws.Range("A23").Select()
